# "stable code with long term balance"
#
# On Sheet2:
#  - F1 held the formula =1.1 (a fitted time-constant). Replace it with the
#    hard-coded, re-tuned value 0.93 (no formula).
#  - Column E (rows 2-92) fits the motor speed curve with
#    $G$1*(1-EXP(-D#/$F$1))  -- F1 used as a divisor. Re-balance the fit so
#    F1 is used as a multiplier instead:
#    $G$1*(1-EXP(-D#*$F$1))

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# F1: replace the =1.1 formula with the literal value 0.93
$ws.Range("F1").Value = 0.93

# E2:E92 - re-enter each row's formula (relative D-reference) so the engine
# keeps them live and re-groups the contiguous runs as shared formulas, same
# as Excel does when a formula is typed/filled across a column.
for ($r = 2; $r -le 92; $r++) {
    $ws.Range("E$r").Formula = "=`$G`$1*(1-EXP(-D$r*`$F`$1))"
}

$excel.CalculateFullRebuild()
